$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.680.30"
$ws.Range("E2").Value = "  +0.25%  "

# Row 3
$ws.Range("D3").Value = "1.720.46"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9971"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9973"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4925"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.89%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06218"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.74%  "

# Row 10
$ws.Range("D10").Value = "1.731.09"
$ws.Range("E10").Value = "  -0.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6096"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.502"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.29%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9972"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "

# Row 17
$ws.Range("D17").Value = "26.490.70"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("E18").Value = "  -0.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007187"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.32%  "

# Row 21
$ws.Range("D21").Value = "1.943.01"
$ws.Range("E21").Value = "  -1.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.434"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.543"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.091"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.80%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.84%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.403"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.742"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.921"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.70%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07972"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.63%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.663"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04500"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.09%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.27%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6251"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.87%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9392"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.78%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.002"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.27%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9968"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.62%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01509"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.574"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.92%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.24%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3852"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "

# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.907"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1158"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.54%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05377"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.764"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.36%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.235"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.30%  "
